# CryCompanywiseStockReport_1: correct stock-report figures.
# - Rows 10/11: the two "SIG-3W/3w Lilliput LED Torch & Table Lamp" line
#   items had their Code/Name/Rate/MRP/Qty/Value fields swapped; restore
#   each row's own data.
# - Various item rows: quantity (F) and value (G = Rate * Qty) corrected.
# - "Sub Total:" / "Grand Total:" rows (B column): recomputed to match the
#   corrected item values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 47438
$ws.Range("C10").Value = 'SIG-3w Lilliput LED Torch &amp; Table Lamp'
$ws.Range("D10").Value = 401.81
$ws.Range("E10").Value = 480.05
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 803.62

$ws.Range("B11").Value = 59408
$ws.Range("C11").Value = 'SIG-3W Lilliput LED Torch &amp; Table Lamp'
$ws.Range("D11").Value = 388.17
$ws.Range("E11").Value = 463.78
$ws.Range("F11").Value = 9
$ws.Range("G11").Value = 3493.53

$ws.Range("F27").Value = 61
$ws.Range("G27").Value = 1562.21

$ws.Range("F35").Value = 129
$ws.Range("G35").Value = 6608.67

$ws.Range("B46").Value = 26768.05

$ws.Range("F48").Value = 29
$ws.Range("G48").Value = 5706.33

$ws.Range("F55").Value = 211
$ws.Range("G55").Value = 40699.79

$ws.Range("F56").Value = 80
$ws.Range("G56").Value = 3369.6

$ws.Range("F58").Value = 37
$ws.Range("G58").Value = 1306.84

$ws.Range("F61").Value = 48
$ws.Range("G61").Value = 1213.92

$ws.Range("F63").Value = 26
$ws.Range("G63").Value = 976.04

$ws.Range("F68").Value = 281
$ws.Range("G68").Value = 26284.74

$ws.Range("B85").Value = 155363.45

$ws.Range("F136").Value = 54
$ws.Range("G136").Value = 1944

$ws.Range("F147").Value = 40
$ws.Range("G147").Value = 4049.2

$ws.Range("F149").Value = 10
$ws.Range("G149").Value = 154.2

$ws.Range("F150").Value = 19
$ws.Range("G150").Value = 586.34

$ws.Range("F157").Value = 25
$ws.Range("G157").Value = 803.5

$ws.Range("B159").Value = 70661.75999999999

$ws.Range("F164").Value = 26
$ws.Range("G164").Value = 1599.52

$ws.Range("B166").Value = 1927.87

$ws.Range("F174").Value = 234
$ws.Range("G174").Value = 4605.12

$ws.Range("F179").Value = 72
$ws.Range("G179").Value = 3205.44

$ws.Range("B180").Value = 35929.88

$ws.Range("F184").Value = 73
$ws.Range("G184").Value = 8762.92

$ws.Range("F197").Value = 16
$ws.Range("G197").Value = 1672.8

$ws.Range("B198").Value = 43649.1

$ws.Range("F257").Value = 0
$ws.Range("G257").Value = 0

$ws.Range("F266").Value = 34
$ws.Range("G266").Value = 1424.94

$ws.Range("B267").Value = 20932.74

$ws.Range("F287").Value = 15
$ws.Range("G287").Value = 3461.4

$ws.Range("F288").Value = 126
$ws.Range("G288").Value = 3119.76

$ws.Range("F290").Value = 88
$ws.Range("G290").Value = 16837.04

$ws.Range("B298").Value = 108730.64

$ws.Range("F324").Value = 43
$ws.Range("G324").Value = 3368.19

$ws.Range("F344").Value = 24
$ws.Range("G344").Value = 2662.56

$ws.Range("B349").Value = 140414.28

$ws.Range("F355").Value = 36
$ws.Range("G355").Value = 6220.44

$ws.Range("F361").Value = 51
$ws.Range("G361").Value = 3759.21

$ws.Range("F363").Value = 36
$ws.Range("G363").Value = 2500.92

$ws.Range("F399").Value = 113
$ws.Range("G399").Value = 6633.1

$ws.Range("F409").Value = 163
$ws.Range("G409").Value = 27926.79

$ws.Range("F414").Value = 16
$ws.Range("G414").Value = 2289.6

$ws.Range("F420").Value = 78
$ws.Range("G420").Value = 7225.14

$ws.Range("B423").Value = 141587.69

$ws.Range("F426").Value = 19
$ws.Range("G426").Value = 3488.02

$ws.Range("B437").Value = 19915.78

$ws.Range("F486").Value = 93
$ws.Range("G486").Value = 5645.1

$ws.Range("B497").Value = 36549.97

$ws.Range("F526").Value = 754
$ws.Range("G526").Value = 72836.39999999999

$ws.Range("F527").Value = 169
$ws.Range("G527").Value = 6288.49

$ws.Range("F528").Value = 191
$ws.Range("G528").Value = 4696.69

$ws.Range("B532").Value = 151507.17

$ws.Range("F549").Value = 32
$ws.Range("G549").Value = 399.04

$ws.Range("B556").Value = 14508.81

$ws.Range("F608").Value = 25
$ws.Range("G608").Value = 1209

$ws.Range("F609").Value = 32
$ws.Range("G609").Value = 3143.36

$ws.Range("B610").Value = 55921.25

$ws.Range("F621").Value = 235
$ws.Range("G621").Value = 14269.2

$ws.Range("F625").Value = 23
$ws.Range("G625").Value = 1479.36

$ws.Range("F627").Value = 14
$ws.Range("G627").Value = 343

$ws.Range("B638").Value = 146107.09

$ws.Range("F640").Value = 4
$ws.Range("G640").Value = 213.76

$ws.Range("B644").Value = 420.81

$ws.Range("F663").Value = 64
$ws.Range("G663").Value = 5043.2

$ws.Range("F664").Value = 16
$ws.Range("G664").Value = 1208.96

$ws.Range("B667").Value = 25599.25

$ws.Range("F672").Value = 68
$ws.Range("G672").Value = 18094.12

$ws.Range("F674").Value = 96
$ws.Range("G674").Value = 5009.28

$ws.Range("B688").Value = 86589.66

$ws.Range("F713").Value = 24
$ws.Range("G713").Value = 3133.2

$ws.Range("F717").Value = 38
$ws.Range("G717").Value = 1033.6

$ws.Range("F718").Value = 112
$ws.Range("G718").Value = 3046.4

$ws.Range("F719").Value = 100
$ws.Range("G719").Value = 2720

$ws.Range("B720").Value = 28518.64

$ws.Range("F758").Value = 10
$ws.Range("G758").Value = 8289.200000000001

$ws.Range("B773").Value = 134345.28

$ws.Range("F778").Value = 10
$ws.Range("G778").Value = 715.6

$ws.Range("F780").Value = 105
$ws.Range("G780").Value = 8977.5

$ws.Range("F781").Value = 9
$ws.Range("G781").Value = 548.37

$ws.Range("B785").Value = 15107.55

$ws.Range("F813").Value = 55
$ws.Range("G813").Value = 5529.15

$ws.Range("F814").Value = 87
$ws.Range("G814").Value = 12208.71

$ws.Range("B815").Value = 37946.64

$ws.Range("F824").Value = 60
$ws.Range("G824").Value = 12953.4

$ws.Range("F833").Value = 135
$ws.Range("G833").Value = 14962.05

$ws.Range("B837").Value = 191044.68

$ws.Range("F843").Value = 65
$ws.Range("G843").Value = 7072.65

$ws.Range("F844").Value = 13
$ws.Range("G844").Value = 329.29

$ws.Range("F849").Value = 32
$ws.Range("G849").Value = 4682.56

$ws.Range("F865").Value = 90
$ws.Range("G865").Value = 4490.1

$ws.Range("F866").Value = 48
$ws.Range("G866").Value = 2735.52

$ws.Range("B867").Value = 200910.41

$ws.Range("B923").Value = 2563612

$ws.Range("B924").Value = 2563612
